$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @"
questions = [
    {
        "title": "In a custom NetSuite sales order form, the \"Payment Terms\" field offers options like \"Net 30,\" \"Net 60,\" and \"COD.\" The \"Credit Limit\" field should be mandatory when \"Payment Terms\" is set to \"Net 60.\" However, sales orders are sometimes saved without filling in the \"Credit Limit\" field, causing incomplete data. Which action should you take to ensure that the \"Credit Limit\" field is always filled when \"Payment Terms\" is \"Net 60\" in the sales order form?",
        "ques_type": 2,
        "options": [
            "Apply a client script to validate the \"Credit Limit\" field and display an error message if it is empty when \"Payment Terms\" is \"Net 60.\"",
            "Create a custom script to hide the \"Credit Limit\" field when \"Payment Terms\" is set to anything other than \"Net 60.\"",
            "Use a workflow to automatically populate the \"Credit Limit\" field with a default value when \"Payment Terms\" is set to \"Net 60.\"",
            "Set up field-level permission that restricts the sales team from saving the form if the \"Credit Limit\" field is empty and \"Payment Terms\" is \"Net 60.\""
        ],
        "score": "Apply a client script to validate the \"Credit Limit\" field and display an error message if it is empty when \"Payment Terms\" is \"Net 60.\""
    },
    {
        "title": "You are a NetSuite administrator responsible for importing complex data sets into NetSuite. Your company recently acquired a new subsidiary, and you must import various record types. During the data import process, you discover that some customer records in the legacy system have duplicate email addresses.Which approach would be most suitable to handle this issue?",
        "ques_type": 2,
        "options": [
            "Manually review and merge duplicate customer records after the data import is complete.",
            "Modify the mapping settings to exclude customer records with duplicate email addresses from the data import.",
            "Use the NetSuite Data Import Assistant to identify and merge duplicate customer records during the import process.",
            "Create a custom script to automatically identify and merge duplicate customer records based on email addresses during the data import."
        ],
        "score": "Use the NetSuite Data Import Assistant to identify and merge duplicate customer records during the import process."
    },
    {
        "title": "You are the NetSuite administrator for a multinational company. The CFO has requested a specific budget report comparing the actual expenses and the budgeted amounts for each department. The company recently underwent a reorganization, and some departments were merged while others were split into separate entities. As a result, the budget structure and department codes have been updated to reflect the changes.Which of the following steps should you take?",
        "ques_type": 2,
        "options": [
            "Create a custom financial statement with a column layout showing each department's actual expenses and budgeted amounts.",
            "Generate a standard budget versus actual report and manually adjust the department codes to match the updated structure.",
            "Run a variance report between the actual expenses and budgeted amounts without considering the department code changes.",
            "Modify the saved search criteria to include the updated department codes and ensure data consistency."
        ],
        "score": "Modify the saved search criteria to include the updated department codes and ensure data consistency."
    },
    {
        "title": "You are an administrator in a company using NetSuite for customer account management. You have been requested to configure the system to allow customers to pay through multiple payment methods and ensure accurate billing. A customer is having trouble adding a new credit card as a payment method in their NetSuite account, despite confirming its validity with their bank.What could be a possible reason for the customer's inability to add the new credit card as a payment method in their account?",
        "ques_type": 2,
        "options": [
            "The customer's browser version is outdated and incompatible with NetSuite's payment gateway.",
            "The customer's account is not set up for online payments.",
            "The customer has reached their maximum limit for the number of payment methods allowed.",
            "The credit card type is not supported by NetSuite's payment gateway integration."
        ],
        "score": "The credit card type is not supported by NetSuite's payment gateway integration."
    }
]
"@

# Remove the trailing newline the here-string adds
$newText = $newText.TrimEnd("`r", "`n")

# Clear the old row 2 (shared string cell) entirely, and reset row1 formatting
$ws.Range("A2").ClearContents()
$ws.Range("A1").ClearFormats()

$cell = $ws.Range("A1")
$cell.Value = $newText
$ws.Rows(1).AutoFit()
